$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the "Valor Mora" (total) figure and the counters.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 973027      # VALOR MORA total
$ws.Range("C13").Value = 10          # Cant. Trabajadores
$ws.Range("F13").Value = 3           # Cant. Periodos

# ------------------------------------------------------------------
# 2. Make room for two extra detail rows under the existing block
#    (rows 16-24 -> 16-26) while preserving the per-row formatting:
#      - rows 16-23 : "normal" data-row style
#      - row 24     : was the closing ("thick border") row, becomes
#                      a normal row
#      - row 25     : new normal row
#      - row 26     : new closing row (keeps the thick-border style
#                      that row 24 used to have)
#    and shift the signature block (previously rows 29-30) down to
#    rows 31-32.
# ------------------------------------------------------------------
$ws.Range("B25:J26").Insert() | Out-Null

# Preserve the "closing row" formatting (currently still on row 24)
# by copying it down onto the new last row (26) first.
$ws.Range("B24:J24").Copy() | Out-Null
$ws.Range("B26:J26").PasteSpecial(-4122) | Out-Null

# Turn row 24 (and the new row 25) into normal data rows by copying
# the formatting from row 23 (a normal row).
$ws.Range("B23:J23").Copy() | Out-Null
$ws.Range("B24:J24").PasteSpecial(-4122) | Out-Null
$ws.Range("B25:J25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 3. Re-write the worker detail rows (16-24) with the new dataset.
# ------------------------------------------------------------------
$ws.Range("C16").Value = "72071604"
$ws.Range("D16").Value = "WILFRIDO UTRIA BERDUGO"
$ws.Range("F16").Value = 35467
$ws.Range("G16").Value = 1400000

$ws.Range("C17").Value = "22639968"
$ws.Range("D17").Value = "DIVINA ESTHER OLMOS MARTINEZ"
$ws.Range("F17").Value = 19000
$ws.Range("G17").Value = 750000

$ws.Range("C18").Value = "72072531"
$ws.Range("D18").Value = "ANTONIO JOSE PUENTES RUIZ"
$ws.Range("F18").Value = 17480
$ws.Range("G18").Value = 690000

$ws.Range("C19").Value = "30685504"
$ws.Range("D19").Value = "LUZ MILEIDA LEON PACHECO"
$ws.Range("F19").Value = 17480
$ws.Range("G19").Value = 690000

$ws.Range("C20").Value = "8641284"
$ws.Range("D20").Value = "LUIS EDUARDO OLMOS MARTINEZ"
$ws.Range("F20").Value = 17480
$ws.Range("G20").Value = 690000

$ws.Range("C21").Value = "1045231139"
$ws.Range("D21").Value = "ARIEL UTRIA CORTINA"
$ws.Range("F21").Value = 17480
$ws.Range("G21").Value = 690000

$ws.Range("C22").Value = "8638995"
$ws.Range("D22").Value = "DIGNO ANTONIO OLMOS MARTINEZ"
$ws.Range("F22").Value = 17480
$ws.Range("G22").Value = 690000

$ws.Range("C23").Value = "3777781"
$ws.Range("D23").Value = "URCINO ANGULO CASTRO"
$ws.Range("F23").Value = 17480
$ws.Range("G23").Value = 690000

$ws.Range("C24").Value = "8527520"
$ws.Range("D24").Value = "LUIS ALFREDO MERCADO CASTELLAR"
$ws.Range("F24").Value = 17480
$ws.Range("G24").Value = 690000

# Column B (Tipo Doc Trabajador) & column E (Periodo Mora) stay "CC" / "1701"
# for every one of those rows - already so from the copy/paste above, but
# make sure explicitly.
foreach ($r in 16..24) {
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 5).Value = "1701"
}

# ------------------------------------------------------------------
# 4. Two new company-level ("NIT") rows (25 & 26).
# ------------------------------------------------------------------
$ws.Range("B25").Value = "NIT"
$ws.Range("C25").Value = "8000166563"
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = "1606"
$ws.Range("F25").Value = 398100
$ws.Range("G25").Value = 0

$ws.Range("B26").Value = "NIT"
$ws.Range("C26").Value = "8000166563"
$ws.Range("D26").Value = ""
$ws.Range("E26").Value = "1605"
$ws.Range("F26").Value = 398100
$ws.Range("G26").Value = 0

Write-Output "edit complete"
